$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Wei-jie%Guan%NULL%0, Zheng-yi%Ni%NULL%0, Zheng-yi%Ni%NULL%0, Yu%Hu%NULL%0, Wen-hua%Liang%NULL%0, Chun-quan%Ou%NULL%0, Jian-xing%He%NULL%0, Lei%Liu%NULL%0, Hong%Shan%NULL%0, Chun-liang%Lei%NULL%0, David S.C.%Hui%NULL%0, Bin%Du%NULL%0, Lan-juan%Li%NULL%0, Guang%Zeng%NULL%0, Kwok-Yung%Yuen%NULL%0, Ru-chong%Chen%NULL%0, Chun-li%Tang%NULL%0, Tao%Wang%NULL%0, Ping-yan%Chen%NULL%0, Jie%Xiang%NULL%0, Shi-yue%Li%NULL%0, Jin-lin%Wang%NULL%0, Zi-jing%Liang%NULL%0, Yi-xiang%Peng%NULL%0, Li%Wei%NULL%0, Yong%Liu%NULL%0, Ya-hua%Hu%NULL%0, Peng%Peng%NULL%0, Jian-ming%Wang%NULL%0, Ji-yang%Liu%NULL%0, Zhong%Chen%NULL%0, Gang%Li%NULL%0, Zhi-jian%Zheng%NULL%0, Shao-qin%Qiu%NULL%0, Jie%Luo%NULL%0, Chang-jiang%Ye%NULL%0, Shao-yong%Zhu%NULL%0, Nan-shan%Zhong%NULL%0]'
$ws.Range("I2").Value = ''
$ws.Range("J2").Value = 'Massachusetts Medical Society'

$ws.Range("E3").Value = '[Ping%Wu%NULL%0, Fang%Duan%NULL%1, Chunhua%Luo%NULL%0, Qiang%Liu%NULL%0, Xingguang%Qu%NULL%1, Liang%Liang%NULL%1, Kaili%Wu%NULL%1]'
$ws.Range("I3").Value = ''
$ws.Range("J3").Value = 'American Medical Association'

$ws.Range("D4").Value = 'This is a cross-sectional study of patients who received a COVID-19 diagnosis between December 30, 2019 and February 7, 2020 at Tongji Hospital.
 A total of 102 patients (48 Male [47%] and 54 Female [53%]) with clinical symptoms, Rt, and chest Computed Tomography (CT) abnormalities were identified with a clinical diagnosis of COVID-19. Patients had a mean [SD] gestational age of 57.63 [14.90] years.
 Of a total of 102 patients identified, 72 patients (36 men [50%] and 36 women [50%]; mean [SD] age, 58.68 [14.81] years) were confirmed to have COVID-19 by laboratory diagnosis with a SARS-CoV-2 RT-PCR assay.
 Only two patients (2.78%) with conjunctivitis were identified from 72 patients with a laboratory confirmed COVID-19. Of those two patients, SARS-CoV-2 RNA fragments were found in ocular discharges by SARS-CoV-2 RT-PCR in only one patient.
 Our findings suspect the incidence of SARS-CoV-2 infection through the ocular surface is extremely low, while the nosocomial infection of SARS-CoV-2 through the eyes after occupational exposure is a potential route.
 To lower the SARS-CoV-2 nosocomial infection, all health care professionals should wear protective goggles.
 The inefficient diagnostic method and the sampling time lag may contribute to the lower positive rate of conjunctival swab samples of SARS-CoV-2.'
$ws.Range("E4").Value = '[Xian%Zhang%NULL%0, Xuhui%Chen%NULL%3, Liwen%Chen%NULL%0, Chaohua%Deng%NULL%3, Xiaojing%Zou%NULL%0, Xiaojing%Zou%NULL%0, Weiyong%Liu%NULL%3, Weiyong%Liu%NULL%0, Huimin%Yu%NULL%3, Huimin%Yu%NULL%0, Bo%Chen%NULL%0, Xufang%Sun%NULL%2]'
$ws.Range("I4").Value = ''
$ws.Range("J4").Value = 'Elsevier Inc.'

$ws.Range("E5").Value = '[Jianhua%Xia%NULL%0, Jianping%Tong%NULL%1, Mengyun%Liu%NULL%1, Ye%Shen%idrshen@zju.edu.cn%1, Dongyu%Guo%guodongyu666@sina.com%1]'
$ws.Range("I5").Value = ''
$ws.Range("J5").Value = 'John Wiley and Sons Inc.'

$ws.Range("D6").Value = 'id="Par1">The 2019 novel coronavirus disease (COVID-19) caused by severe acute respiratory syndrome coronavirus 2 (SARS-CoV-2) has spread globally, while the routes of transmission of this virus are still controversial.
 We enrolled 33 patients, without any ocular manifestation, with their ocular surface swabs collected for virus detection.
 RNA was detected strong positive in samples of both eyes from two patients.
 Therefore, SARS-CoV-2 may exist in the normal ocular surface of COVID-19 patients, suggesting that this virus might be spread through conjunctival contact.
'
$ws.Range("E6").Value = '[Hua-Tao%Xie%NULL%0, Shi-Yun%Jiang%NULL%1, Kang-Kang%Xu%NULL%1, Xin%Liu%NULL%0, Bing%Xu%NULL%1, Lin%Wang%lin_wang@hust.edu.cn%0, Ming-Chang%Zhang%mingchangzhang@hotmail.com%1]'
$ws.Range("I6").Value = ''
$ws.Range("J6").Value = 'BioMed Central'

$ws.Range("E7").Value = '[Liwen%Chen%NULL%0, Chaohua%Deng%NULL%0, Chaohua%Deng%NULL%0, Xuhui%Chen%NULL%0, Xuhui%Chen%NULL%0, Xian%Zhang%NULL%0, Bo%Chen%NULL%0, Huimin%Yu%NULL%0, Yuanjun%Qin%NULL%1, Ke%Xiao%NULL%1, Hong%Zhang%tjyksys@163.com%1, Xufang%Sun%sunxufang2016@163.com%0]'
$ws.Range("I7").Value = ''
$ws.Range("J7").Value = 'John Wiley and Sons Inc.'

$ws.Range("E8").Value = '[Kiran%Kumar%NULL%0, Akshata A%Prakash%NULL%1, Suresh Babu%Gangasagara%NULL%1, Sujatha BL%Rathod%NULL%1, K%Ravi%NULL%1, Ambica%Rangaiah%NULL%1, Sathyanarayan Muthur%Shankar%NULL%1, Shantala Gowdara%Basawarajappa%NULL%1, Shashi%Bhushan%NULL%1, T G%Neeraja%NULL%1, Srinivas%Khandenahalli%NULL%1, M%Swetha%NULL%1, Priyam%Gupta%NULL%1, U C%Sampritha%NULL%1, Guru N S%Prasad%NULL%1, Chakravarthy Raghunathan%Jayanthi%NULL%1]'
$ws.Range("I8").Value = ''
$ws.Range("J8").Value = 'Wolters Kluwer - Medknow'

$ws.Range("E9").Value = '[Yunyun%Zhou%NULL%0, Chaoye%Duan%NULL%1, Yuyang%Zeng%NULL%1, Yongqing%Tong%NULL%1, Yuhong%Nie%NULL%1, Yang%Yang%NULL%0, Zhen%Chen%NULL%1, Changzheng%Chen%chenchangzheng@whu.edu.cn%1]'
$ws.Range("I9").Value = ''
$ws.Range("J9").Value = 'by the American Academy of Ophthalmology'

$ws.Range("D10").Value = 'Objectives
id="Par1">Since there are few reports on the ocular involvement of coronavirus disease 2019 (COVID-19) patients, this study aimed to assess the presence of severe acute respiratory syndrome Coronavirus-2 (SARS-CoV-2) in the tears of patients with COVID-19.
Methods
id="Par2">In this prospective case series, nasopharyngeal and tear sampling of 43 patients with severe COVID-19 were performed.

 The quantitative reverse transcription polymerase chain reaction (RT-PCR) was conducted to detect SARS-CoV-2. Ocular and systemic signs and symptoms were recorded from their medical history.


Results
id="Par3">The mean age of patients was 56 ± 13 years.

 The average disease time from initiation of symptoms was 3.27 days, range: 1–7 days.

 Forty-one patients (95.3%) had fever at the time of sampling.

 Only one patient had conjunctivitis.

 Thirty (69.8%) nasopharyngeal and three (7%) tear samples were positive for SARS-CoV-2. The result of tear sample was positive in the patient with conjunctivitis.

 All patients with positive tear RT-PCR results had positive nasopharyngeal RT-PCR results.


Conclusions
id="Par4">Ocular manifestation was rare in this series of severe COVID-19 patients, however, 7% of the patients had viral RNA in their conjunctival secretions.

 Therefore, possibility of ocular transmission should be considered even in the absence of ocular manifestations.


'
$ws.Range("E10").Value = '[Saeed%Karimi%NULL%0, Amir%Arabi%amir_arab_91@yahoo.com%1, Toktam%Shahraki%NULL%2, Toktam%Shahraki%NULL%0, Sare%Safi%NULL%1]'
$ws.Range("I10").Value = ''
$ws.Range("J10").Value = 'Nature Publishing Group UK'

$ws.Range("D11").Value = 'Background
id="Par1">A recent increase in children admitted with hypotensive shock and fever in the context of the COVID-19 outbreak requires an urgent characterization and assessment of the involvement of SARS-CoV-2 infection.

 This is a case series performed at 4 academic tertiary care centers in Paris of all the children admitted to the pediatric intensive care unit (PICU) with shock, fever and suspected SARS-CoV-2 infection between April 15th and April 27th, 2020.
Results
id="Par2">20 critically ill children admitted for shock had an acute myocarditis (left ventricular ejection fraction, 35% (25–55); troponin, 269 ng/mL (31–4607)), and arterial hypotension with mainly vasoplegic clinical presentation.

 The first symptoms before PICU admission were intense abdominal pain and fever for 6 days (1–10).

 All children had highly elevated C-reactive protein (&gt; 94 mg/L) and procalcitonin (&gt; 1.6 ng/mL) without microbial cause.

 At least one feature of Kawasaki disease was found in all children (fever, n = 20, skin rash, n = 10; conjunctivitis, n = 6; cheilitis, n = 5; adenitis, n = 2), but none had the typical form.

 SARS-CoV-2 PCR and serology were positive for 10 and 15 children, respectively.

 One child had both negative SARS-CoV-2 PCR and serology, but had a typical SARS-CoV-2 chest tomography scan.

 All children but one needed an inotropic/vasoactive drug support (epinephrine, n = 12; milrinone, n = 10; dobutamine, n = 6, norepinephrine, n = 4) and 8 were intubated.

 All children received intravenous immunoglobulin (2 g per kilogram) with adjuvant corticosteroids (n = 2), IL 1 receptor antagonist (n = 1) or a monoclonal antibody against IL-6 receptor (n = 1).

 All children survived and were afebrile with a full left ventricular function recovery at PICU discharge.


Conclusions
id="Par3">Acute myocarditis with intense systemic inflammation and atypical Kawasaki disease is an emerging severe pediatric disease following SARS-CoV-2 infection.

 Early recognition of this disease is needed and referral to an expert center is recommended.

 A delayed and inappropriate host immunological response is suspected.

 While underlying mechanisms remain unclear, further investigations are required to target an optimal treatment.


'
$ws.Range("E11").Value = '[Marion%Grimaud%NULL%0, Julie%Starck%NULL%2, Michael%Levy%NULL%0, Clémence%Marais%NULL%2, Judith%Chareyre%NULL%2, Diala%Khraiche%NULL%2, Marianne%Leruez-Ville%NULL%2, Pierre%Quartier%NULL%2, Pierre Louis%Léger%NULL%2, Guillaume%Geslain%NULL%2, Nada%Semaan%NULL%2, Florence%Moulin%NULL%2, Matthieu%Bendavid%NULL%2, Sandrine%Jean%NULL%2, Géraldine%Poncelet%NULL%2, Sylvain%Renolleau%NULL%2, Mehdi%Oualha%mehdi.oualha@aphp.fr%3]'
$ws.Range("I11").Value = ''
$ws.Range("J11").Value = 'Springer International Publishing'

$ws.Range("E12").Value = '[Zhixiong%Fang%NULL%0, Yi%Zhang%NULL%0, Changfa%Hang%NULL%1, Jingwen%Ai%NULL%1, Shaojie%Li%NULL%1, Wenhong%Zhang%NULL%0]'
$ws.Range("I12").Value = ''
$ws.Range("J12").Value = 'The British Infection Association. Published by Elsevier Ltd.'

$ws.Range("E13").Value = '[Paola%Valente%NULL%0, Giancarlo%Iarossi%NULL%1, Matteo%Federici%NULL%1, Sergio%Petroni%NULL%1, Paolo%Palma%NULL%1, Nicola%Cotugno%NULL%1, Maria A.%De Ioris%NULL%1, Andrea%Campana%NULL%1, Luca%Buzzonetti%NULL%1]'
$ws.Range("I13").Value = ''
$ws.Range("J13").Value = 'American Association for Pediatric Ophthalmology and Strabismus. Published by Elsevier Inc.'

$ws.Range("E14").Value = '[Ivan Yu Jun%Seah%NULL%0, Danielle E.%Anderson%NULL%1, Adrian Eng Zheng%Kang%NULL%1, Linfa%Wang%NULL%1, Pooja%Rao%NULL%1, Barnaby Edward%Young%NULL%0, David Chien%Lye%NULL%0, Rupesh%Agrawal%Rupesh_agrawal@ttsh.com.sg%1]'
$ws.Range("I14").Value = ''
$ws.Range("J14").Value = 'by the American Academy of Ophthalmology'
